# Casper Tak's "First order list" - regulator chip swap
# Removed old regulator circuits; a better (tested) chip found: RT7272A.
# Add a new "stepper motor protector" line item to the order list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the Link column first so the URL string lands in the shared-string
# table ahead of the component name, matching save order.
$ws.Range("D9").Value = "https://nl.aliexpress.com/item/32917996858.html?gatewayAdapt=glo2nld"
$ws.Range("B9").Value = "stepper motor protector"
$ws.Range("C9").Value = 6

# Leave the selection on A9, as it was when the workbook was last saved.
$ws.Range("A9").Select()
